$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.491026518625059
$ws.Range("D2").Value = 5.883356914042586
$ws.Range("E2").Value = 12.57342016251977
$ws.Range("F2").Value = 30.55885026450516
$ws.Range("G2").Value = 39.86534104622985
$ws.Range("H2").Value = 16.83316682029627
$ws.Range("I2").Value = 27.61652329687077
$ws.Range("K2").Value = 16.9908014564247
$ws.Range("L2").Value = 9.577868010516578
$ws.Range("N2").Value = 17.65228382869239

$ws.Range("C3").Value = 9.427497200206613
$ws.Range("D3").Value = 5.891914011016291
$ws.Range("E3").Value = 12.52516203091558
$ws.Range("F3").Value = 30.51539782244329
$ws.Range("G3").Value = 39.75378929431074
$ws.Range("H3").Value = 16.87463977238474
$ws.Range("I3").Value = 27.62629622543163
$ws.Range("K3").Value = 16.53518966051216
$ws.Range("L3").Value = 9.581003474927124
$ws.Range("N3").Value = 17.71083465628408

$ws.Range("C4").Value = 9.390112568384614
$ws.Range("D4").Value = 5.897483163426196
$ws.Range("E4").Value = 12.49821140890673
$ws.Range("F4").Value = 30.49869935323659
$ws.Range("G4").Value = 39.70052272482256
$ws.Range("H4").Value = 16.90398337401394
$ws.Range("I4").Value = 27.64016427238827
$ws.Range("K4").Value = 16.25203933559521
$ws.Range("L4").Value = 9.584647709655489
$ws.Range("N4").Value = 17.74872684430441

$ws.Range("C5").Value = 9.375298689625714
$ws.Range("D5").Value = 5.899832164939745
$ws.Range("E5").Value = 12.48791033608076
$ws.Range("F5").Value = 30.49440446930549
$ws.Range("G5").Value = 39.68265146049826
$ws.Range("H5").Value = 16.91691331273641
$ws.Range("I5").Value = 27.64778904459716
$ws.Range("K5").Value = 16.13597413776932
$ws.Range("L5").Value = 9.586565113669341
$ws.Range("N5").Value = 17.76465769411549

$ws.Range("C6").Value = 9.372864619252779
$ws.Range("D6").Value = 5.900227027024139
$ws.Range("E6").Value = 12.48624121893585
$ws.Range("F6").Value = 30.49384284340643
$ws.Range("G6").Value = 39.67991566972643
$ws.Range("H6").Value = 16.91911893875692
$ws.Range("I6").Value = 27.64917415491856
$ws.Range("K6").Value = 16.1166659212726
$ws.Range("L6").Value = 9.586909609217031
$ws.Range("N6").Value = 17.76733259861167

$ws.Range("C7").Value = 9.389911063789473
$ws.Range("D7").Value = 5.897514520497642
$ws.Range("E7").Value = 12.49806971614243
$ws.Range("F7").Value = 30.49863127069931
$ws.Range("G7").Value = 39.70026617360251
$ws.Range("H7").Value = 16.90415381981697
$ws.Range("I7").Value = 27.64025911987317
$ws.Range("K7").Value = 16.25047654953368
$ws.Range("L7").Value = 9.584671817886361
$ws.Range("N7").Value = 17.74893970985793

$ws.Range("C8").Value = 9.468792618350692
$ws.Range("D8").Value = 5.886242190703857
$ws.Range("E8").Value = 12.55622931761461
$ws.Range("F8").Value = 30.54179696892173
$ws.Range("G8").Value = 39.82372327986332
$ws.Range("H8").Value = 16.84665970413937
$ws.Range("I8").Value = 27.61825754056254
$ws.Range("K8").Value = 16.83451573477212
$ws.Range("L8").Value = 9.578592422299421
$ws.Range("N8").Value = 17.67206986919679

$ws.Range("C9").Value = 9.63576858993121
$ws.Range("D9").Value = 5.866623695347062
$ws.Range("E9").Value = 12.69117462581544
$ws.Range("F9").Value = 30.70552559582144
$ws.Range("G9").Value = 40.18609055867936
$ws.Range("H9").Value = 16.76483056890127
$ws.Range("I9").Value = 27.63770188933514
$ws.Range("K9").Value = 17.94555917734699
$ws.Range("L9").Value = 9.580301885294505
$ws.Range("N9").Value = 17.53668065843045

$ws.Range("C10").Value = 9.765115487436047
$ws.Range("D10").Value = 5.853707762482837
$ws.Range("E10").Value = 12.80250938390955
$ws.Range("F10").Value = 30.8736932651946
$ws.Range("G10").Value = 40.52452492395693
$ws.Range("H10").Value = 16.72373760375522
$ws.Range("I10").Value = 27.69028930547693
$ws.Range("K10").Value = 18.73199953253021
$ws.Range("L10").Value = 9.589848279088104
$ws.Range("N10").Value = 17.44649390612975

$ws.Range("C11").Value = 9.825222156676096
$ws.Range("D11").Value = 5.84815346453054
$ws.Range("E11").Value = 12.85567351477894
$ws.Range("F11").Value = 30.96046433789519
$ws.Range("G11").Value = 40.69381132046838
$ws.Range("H11").Value = 16.70920921694625
$ws.Range("I11").Value = 27.722532339133
$ws.Range("K11").Value = 19.08162177838196
$ws.Range("L11").Value = 9.595983848097758
$ws.Range("N11").Value = 17.40746586533653

$ws.Range("C12").Value = 9.848148959530992
$ws.Range("D12").Value = 5.846096104949227
$ws.Range("E12").Value = 12.87615518817705
$ws.Range("F12").Value = 30.99478395757822
$ws.Range("G12").Value = 40.760082166008
$ws.Range("H12").Value = 16.70430905766884
$ws.Range("I12").Value = 27.73593622123943
$ws.Range("K12").Value = 19.21272136260139
$ws.Range("L12").Value = 9.598564072735687
$ws.Range("N12").Value = 17.39297320010021

$ws.Range("C13").Value = 9.843204155172609
$ws.Range("D13").Value = 5.846537155337295
$ws.Range("E13").Value = 12.8717287503597
$ws.Range("F13").Value = 30.98732789887723
$ws.Range("G13").Value = 40.74571390199323
$ws.Range("H13").Value = 16.70533760875864
$ws.Range("I13").Value = 27.73299639852296
$ws.Range("K13").Value = 19.18454618056311
$ws.Range("L13").Value = 9.597996970880704
$ws.Range("N13").Value = 17.39608173172401

$ws.Range("C14").Value = 9.827105131961247
$ws.Range("D14").Value = 5.847983285226447
$ws.Range("E14").Value = 12.8573516217795
$ws.Range("F14").Value = 30.9632586460472
$ws.Range("G14").Value = 40.69922029785087
$ws.Range("H14").Value = 16.70879401094629
$ws.Range("I14").Value = 27.72361117083491
$ws.Range("K14").Value = 19.09243404630465
$ws.Range("L14").Value = 9.596190984048441
$ws.Range("N14").Value = 17.40626781039103

$ws.Range("C15").Value = 9.817265106429863
$ws.Range("D15").Value = 5.848875055724534
$ws.Range("E15").Value = 12.84859036068037
$ws.Range("F15").Value = 30.9487053285403
$ws.Range("G15").Value = 40.6710224312933
$ws.Range("H15").Value = 16.71098955064613
$ws.Range("I15").Value = 27.71801786588369
$ws.Range("K15").Value = 19.03584056368047
$ws.Range("L15").Value = 9.595118177677604
$ws.Range("N15").Value = 17.41254434847769

$ws.Range("C16").Value = 9.761211460559952
$ws.Range("D16").Value = 5.854077190468689
$ws.Range("E16").Value = 12.79908459102389
$ws.Range("F16").Value = 30.86822797095194
$ws.Range("G16").Value = 40.51376695563228
$ws.Range("H16").Value = 16.72477109791929
$ws.Range("I16").Value = 27.6883494126319
$ws.Range("K16").Value = 18.70897627076583
$ws.Range("L16").Value = 9.589483301367482
$ws.Range("N16").Value = 17.44908461018528

$ws.Range("C17").Value = 9.72713709527028
$ws.Range("D17").Value = 5.857350618445394
$ws.Range("E17").Value = 12.76935052471981
$ws.Range("F17").Value = 30.82147804261093
$ws.Range("G17").Value = 40.42119644042849
$ws.Range("H17").Value = 16.73429424441005
$ws.Range("I17").Value = 27.67227869411467
$ws.Range("K17").Value = 18.50627977066607
$ws.Range("L17").Value = 9.586485032630634
$ws.Range("N17").Value = 17.47201201648578

$ws.Range("C18").Value = 9.707658592189272
$ws.Range("D18").Value = 5.859263657919337
$ws.Range("E18").Value = 12.75248600970641
$ws.Range("F18").Value = 30.79555656099291
$ws.Range("G18").Value = 40.36939800154482
$ws.Range("H18").Value = 16.74016360191406
$ws.Range("I18").Value = 27.66381866425296
$ws.Range("K18").Value = 18.38893530606252
$ws.Range("L18").Value = 9.584929331697047
$ws.Range("N18").Value = 17.4853874049918

$ws.Range("C19").Value = 9.701084639187956
$ws.Range("D19").Value = 5.859916583725069
$ws.Range("E19").Value = 12.74681717414297
$ws.Range("F19").Value = 30.78694665225165
$ws.Range("G19").Value = 40.35210933462493
$ws.Range("H19").Value = 16.74221809379543
$ws.Range("I19").Value = 27.66108884553479
$ws.Range("K19").Value = 18.34907816263352
$ws.Range("L19").Value = 9.584431621123985
$ws.Range("N19").Value = 17.48994842640489

$ws.Range("C20").Value = 9.730752050089784
$ws.Range("D20").Value = 5.856999027593055
$ws.Range("E20").Value = 12.7724912504075
$ws.Range("F20").Value = 30.82635458938536
$ws.Range("G20").Value = 40.430901373472
$ws.Range("H20").Value = 16.73323991076984
$ws.Range("I20").Value = 27.67390838108435
$ws.Range("K20").Value = 18.52793663478473
$ws.Range("L20").Value = 9.586786737476061
$ws.Range("N20").Value = 17.4695518891465

$ws.Range("C21").Value = 9.831829439479248
$ws.Range("D21").Value = 5.847557277260013
$ws.Range("E21").Value = 12.86156514772628
$ws.Range("F21").Value = 30.97028884756993
$ws.Range("G21").Value = 40.71281815240228
$ws.Range("H21").Value = 16.70776244055754
$ws.Range("I21").Value = 27.72633545697896
$ws.Range("K21").Value = 19.1195256946042
$ws.Range("L21").Value = 9.596714485004561
$ws.Range("N21").Value = 17.40326814663356

$ws.Range("C22").Value = 9.898847513695772
$ws.Range("D22").Value = 5.841654161792867
$ws.Range("E22").Value = 12.92181147777152
$ws.Range("F22").Value = 31.0728666591127
$ws.Range("G22").Value = 40.90966930200562
$ws.Range("H22").Value = 16.69461800029943
$ws.Range("I22").Value = 27.76755785705884
$ws.Range("K22").Value = 19.49856622238605
$ws.Range("L22").Value = 9.604699185727235
$ws.Range("N22").Value = 17.36161686592473

$ws.Range("C23").Value = 9.862996495549082
$ws.Range("D23").Value = 5.844780362873681
$ws.Range("E23").Value = 12.88947526379994
$ws.Range("F23").Value = 31.01734627211694
$ws.Range("G23").Value = 40.8034671349035
$ws.Range("H23").Value = 16.70131182025839
$ws.Range("I23").Value = 27.74492115878796
$ws.Range("K23").Value = 19.29699876211776
$ws.Range("L23").Value = 9.60030105961518
$ws.Range("N23").Value = 17.38369453450821

$ws.Range("C24").Value = 9.729117381096158
$ws.Range("D24").Value = 5.857157884903623
$ws.Range("E24").Value = 12.77107061100033
$ws.Range("F24").Value = 30.82414692459832
$ws.Range("G24").Value = 40.42650934227894
$ws.Range("H24").Value = 16.73371534674715
$ws.Range("I24").Value = 27.67316917218137
$ws.Range("K24").Value = 18.51814808465963
$ws.Range("L24").Value = 9.586649813213967
$ws.Range("N24").Value = 17.47066350797549

$ws.Range("C25").Value = 9.589364369190893
$ws.Range("D25").Value = 5.871666712284545
$ws.Range("E25").Value = 12.65248492975291
$ws.Range("F25").Value = 30.6527899744322
$ws.Range("G25").Value = 40.07527796442373
$ws.Range("H25").Value = 16.78363837230005
$ws.Range("I25").Value = 27.62572062347501
$ws.Range("K25").Value = 17.64962836078432
$ws.Range("L25").Value = 9.578381454836716
$ws.Range("N25").Value = 17.57167137346439

